$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep only the "Who is CEO of SpaceX?" / "Elon Musk" / "Person" row (originally row 4)
# and drop the other four Q&A rows (originally rows 2, 3, 5 and 6).
$ws.Range("A5:C6").EntireRow.Delete()
$ws.Range("A2:C3").EntireRow.Delete()

# Shrink the columns now that the remaining text is much shorter
$ws.Columns("A").ColumnWidth = 20.5
$ws.Columns("B").ColumnWidth = 9.15

# Update the selected cell
$ws.Range("N17").Select()

# Configure page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
